$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 2: update the resource record -------------------------------------
$ws.Range("A2").Value = "Aca_CURCE"
$ws.Range("B2").Value = "Center for Undergraduate Research and Creative Engagement (CURCE)"
$ws.Range("C2").Value = "Acadata()"
$ws.Range("D2").Value = "Center for Undergraduate Research and Creative Engagement (CURCE)"
$ws.Range("E2").Value = "F"
$ws.Range("F2").Value = "T"

# B2/D2 previously carried a special purple, wrap-text style; the new entry
# drops that formatting back to the plain style used elsewhere in column A.
$ws.Range("A2").Copy()
$ws.Range("B2").PasteSpecial(-4122)
$ws.Range("D2").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Recalculate formulas in B4/B5/B6 before snapshotting their text -------
$excel.Calculate()

# --- Rows 8-10: drop in the freshly generated code snippets -----------------
$ws.Range("A8").Value = $ws.Range("B4").Value2
$ws.Range("A9").Value = $ws.Range("B5").Value2
$ws.Range("A10").Value = $ws.Range("B6").Value2

# --- Rows 12-14: clear the old snippet text left over in column A ----------
$ws.Range("A12").Clear()
$ws.Range("A13").Clear()
$ws.Range("A14").Clear()

# --- Selection now rests on A13, matching the post-edit cursor position ----
$ws.Range("A13").Select()

$excel.Calculate()
